# Update the "想去人数" (want-to-go count) figures in column F across the
# three data sheets (展览, 演出, 全部类型) to reflect a refreshed scrape.
# 本地生活 has no changed rows and is left untouched.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 -------------------------------------------------------
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value  = 819
$ws.Range("F3").Value  = 195
$ws.Range("F4").Value  = 1400
$ws.Range("F5").Value  = 834
$ws.Range("F6").Value  = 473
$ws.Range("F7").Value  = 626
$ws.Range("F8").Value  = 188
$ws.Range("F10").Value = 46
$ws.Range("F12").Value = 114
$ws.Range("F13").Value = 1589
$ws.Range("F14").Value = 201
$ws.Range("F15").Value = 31
$ws.Range("F16").Value = 477
$ws.Range("F18").Value = 393
$ws.Range("F19").Value = 108
$ws.Range("F20").Value = 634
$ws.Range("F21").Value = 29
$ws.Range("F22").Value = 215
$ws.Range("F23").Value = 728
$ws.Range("F25").Value = 1445
$ws.Range("F26").Value = 168

# --- Sheet 2: 演出 -------------------------------------------------------
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 651
$ws.Range("F5").Value = 199
$ws.Range("F6").Value = 14
$ws.Range("F8").Value = 64
$ws.Range("F9").Value = 37

# --- Sheet 4: 全部类型 ----------------------------------------------------
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value  = 819
$ws.Range("F4").Value  = 195
$ws.Range("F5").Value  = 1400
$ws.Range("F6").Value  = 834
$ws.Range("F9").Value  = 473
$ws.Range("F10").Value = 626
$ws.Range("F11").Value = 651
$ws.Range("F12").Value = 188
$ws.Range("F14").Value = 46
$ws.Range("F16").Value = 114
$ws.Range("F17").Value = 1589
$ws.Range("F18").Value = 199
$ws.Range("F19").Value = 201
$ws.Range("F20").Value = 31
$ws.Range("F21").Value = 477
$ws.Range("F23").Value = 393
$ws.Range("F24").Value = 14
$ws.Range("F25").Value = 108
$ws.Range("F27").Value = 64
$ws.Range("F28").Value = 634
$ws.Range("F29").Value = 37
$ws.Range("F33").Value = 29
$ws.Range("F34").Value = 215
$ws.Range("F35").Value = 728
$ws.Range("F37").Value = 1445
$ws.Range("F38").Value = 168
